# "funciona todo menos categoria"
#   - Brand: remove the last two rows (6/Nike, 7/emo) that were added previously.
#   - Client: add a new client row (jose / 5552 / dgg).
#   - Category: add a new category row (pesticida / para pesticid).

$wb = $excel.ActiveWorkbook

# --- Brand: drop rows 7 and 8 (the extra Nike/emo rows) ---------------------
$brand = $wb.Worksheets.Item("Brand")
$brand.Range("A7:B8").ClearContents()

# --- Client: append ID=1, Name=jose, DNI=5552, Nickname=dgg -----------------
$client = $wb.Worksheets.Item("Client")
$client.Cells.Item(2, 1).Value = 1
$client.Cells.Item(2, 2).Value = "jose"

# DNI "5552" must be stored as text (not a number) - force text formatting
# before writing the value, then strip the formatting back off so the cell
# keeps the default style while the stored value stays a text string.
$client.Cells.Item(2, 3).NumberFormat = "@"
$client.Cells.Item(2, 3).Value = "5552"
$client.Cells.Item(2, 3).ClearFormats()

$client.Cells.Item(2, 4).Value = "dgg"

# --- Category: append ID=1, Name=pesticida, Description=para pesticid -------
$category = $wb.Worksheets.Item("Category")
$category.Cells.Item(2, 1).Value = 1
$category.Cells.Item(2, 2).Value = "pesticida"
$category.Cells.Item(2, 3).Value = "para pesticid"
